# edit.ps1 - applies the HOPE_OLUWALOLOPE -> FADIPE AL-AMEEN profile update
# described by the commit "FADO added his profile".

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2
    ) | Out-Null
}

# 1. Remove the profile picture anchored in the name/header paragraph.
#    (Shapes.Item(1) resolves to the "Picture 7" drawing; the
#    "Straight Connector 3" line under EDUCATION is Item(2) and stays.)
$d.Shapes.Item(1).Delete()

# 2. Name: "HOPE OLUWALOLOPE" -> "FADIPE AL-AMEEN"
Replace-Text "HOPE" "FADIPE"
Replace-Text " OLUWALOLOPE" " AL-AMEEN"

# 3. Title: "Back-End Web Developer" -> "Front-End Web Developer"
Replace-Text "Back-End Web Developer" "Front-End Web Developer"

# 4. Phone number
Replace-Text ": 08120087057, 09056871401." ": 08023301369."

# 5. Email address (label "Email:" text itself is unchanged)
Replace-Text "blessyn2hope@gmail.com" "fadipetomi00@gmail.com"

# 6. LinkedIn URL
Replace-Text "https://www.linkedin.com/in/hopeoluwalolope/" "https://www.linkedin.com/in/fadipe-al-ameen-a1b51a160/"

# 7. Twitter URL (keep the leading space before the URL)
Replace-Text " https://twitter.com/hopeblessy" " https://twitter.com/Fadipetomi"

# 8. GitHub URL
Replace-Text "https://github.com/ebonyhope" "https://github.com/fadhoo"

# 9. Facebook URL, also swallowing the trailing manual line break "^l"
#    that used to follow it (the new document has no break there).
Replace-Text "hope.oluwalolope^l" "fadipe.tomi"

# 10. Education table updates
Replace-Text "Sciences" "Engineering"
Replace-Text "Computer sciences" "Systems Engineering"
Replace-Text "400" "300"
Replace-Text "Algorithms and Discrete Structures" "Algorithms and Data Structures"
